$d = $word.ActiveDocument

# 1. "CEA201_Test01" -> "CEA201_Test01aaaaaaaaa" (Subject line)
$d.Content.Find.Execute("CEA201_Test01", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CEA201_Test01aaaaaaaaa", 2)

# 2. "Number of question: 2" -> "Number of question: 1"
#    Target only the digit run so the two existing runs
#    ("Number of question: " and "2") stay separate, just like the source
#    document, instead of being coalesced into a single run by a plain
#    Range.Text assignment.
$full = $d.Content.Text
$idx = $full.IndexOf("Number of question: 2")
$charPos = $idx + 20
$numRng = $d.Range($charPos, $charPos + 1)
$xmlPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
          '<pkg:xmlData>' + `
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body><w:p><w:r><w:t>1</w:t></w:r></w:p></w:body>' + `
          '</w:document>' + `
          '</pkg:xmlData></pkg:part></pkg:package>'
$numRng.InsertXML($xmlPkg)

# 3. First question's ANSWER value "D" -> "" (table 1, row 6, col 2)
$table1 = $d.Tables.Item(1)
$table1.Cell(6, 2).Range.Text = ""

# 4. Second question's text replaced with the new Pentium 4 question
$d.Content.Find.Execute("We are test question 2 more powerfull and essily to write", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "The Pentium 4 _________ component executes micro-operations, fetching the required data from the L1 data cache and temporarily storing results in registers.", `
                         2)

# 5. Second question's ANSWER value "C" -> "" (table 2, row 6, col 2)
$table2 = $d.Tables.Item(2)
$table2.Cell(6, 2).Range.Text = ""
